# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gains a new "property_category" column
# (value "stock" for every existing row), inserted between the existing
# "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # "股票" (stock) sheet

# Insert a new column at H, pushing the existing date / legislator_name /
# legislator_id columns one slot to the right (H->I, I->J, J->K).
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Populate every existing data row with the new category value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
